$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 updates
$ws.Range("D16").Value = "image_20250807111728_ppp0.jpg"
$ws.Range("I16").Value = "'642,530,686,574"
$ws.Range("J16").Value = "'0.75"

# Row 17 updates
$ws.Range("D17").Value = "image_20250807111728_ppp0.jpg"
$ws.Range("I17").Value = "'794,481,830,525"
$ws.Range("J17").Value = "'0.70"
